# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.192.94"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.001.53"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0803"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.38%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.297.47"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.848"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "2.005.67"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "37.137.49"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.143"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  +11.70%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0980"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").Value = "1.381.32"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.96%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.46%  "
